$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 9000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H33").Value = 382.6154
$ws.Range("I33").Value = 329.42856
$ws.Range("J33").Value = 444.66666
$ws.Range("K33").Value = 329.42856
$ws.Range("L33").Value = 444.66666
$ws.Range("M33").Value = -100.42856
$ws.Range("N33").Value = -902.66666

$ws.Range("H38").Value = 606.1667
$ws.Range("J38").Value = 1666.6666
$ws.Range("L38").Value = 4999.9998
$ws.Range("N38").Value = -5743.9998

$ws.Range("H40").Value = 7771.778
$ws.Range("I40").Value = 3473
$ws.Range("K40").Value = 3473
$ws.Range("M40").Value = -3298

$ws.Range("H125").Value = 3149.25
$ws.Range("I125").Value = 1187
$ws.Range("K125").Value = 10683
$ws.Range("M125").Value = -8223

$ws.Range("H137").Value = 2321.1875
$ws.Range("I137").Value = 1669.2354
$ws.Range("K137").Value = 5007.706200000001
$ws.Range("M137").Value = -2457.706200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 942.86664
$ws.Range("I2").Value = 295.07144
$ws.Range("K2").Value = 295.07144
$ws.Range("M2").Value = -182.07144

$ws.Range("H7").Value = 50000
$ws.Range("J7").Value = 50000
$ws.Range("L7").Value = 50000
$ws.Range("N7").Value = -50228

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H19").Value = 2333
$ws.Range("J19").Value = 2749.5
$ws.Range("L19").Value = 2749.5
$ws.Range("N19").Value = -3207.5

$ws.Range("H61").Value = 6380.7
$ws.Range("I61").Value = 6312
$ws.Range("K61").Value = 6312
$ws.Range("M61").Value = -6100

$ws.Range("H102").Value = 2258.5833
$ws.Range("I102").Value = 1122.6666
$ws.Range("K102").Value = 1122.6666
$ws.Range("M102").Value = 499.3334

$ws.Range("H110").Value = 1810
$ws.Range("I110").Value = 1782.8572
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1782.8572
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 262.1428000000001
$ws.Range("N110").Value = -6090

$ws.Range("H116").Value = 942.86664
$ws.Range("I116").Value = 295.07144
$ws.Range("K116").Value = 295.07144
$ws.Range("M116").Value = 1998.92856

$ws.Range("H122").Value = 2174.1428
$ws.Range("I122").Value = 2036.5
$ws.Range("K122").Value = 6109.5
$ws.Range("M122").Value = -3659.5

$ws.Range("H132").Value = 2333.4285
$ws.Range("I132").Value = 2464.4
$ws.Range("J132").Value = 2006
$ws.Range("K132").Value = 7393.200000000001
$ws.Range("L132").Value = 6018
$ws.Range("M132").Value = -4863.200000000001
$ws.Range("N132").Value = -11078

$ws.Range("H136").Value = 6380.7
$ws.Range("I136").Value = 6312
$ws.Range("K136").Value = 18936
$ws.Range("M136").Value = -16386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 942.86664
$ws.Range("I3").Value = 295.07144
$ws.Range("K3").Value = 295.07144
$ws.Range("M3").Value = -181.07144

$ws.Range("H81").Value = 13733
$ws.Range("J81").Value = 13733
$ws.Range("L81").Value = 13733
$ws.Range("N81").Value = -15855

$ws.Range("H84").Value = 13733
$ws.Range("J84").Value = 13733
$ws.Range("L84").Value = 41199
$ws.Range("N84").Value = -51807

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2600
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H55").Value = 4573
$ws.Range("I55").Value = 4573
$ws.Range("K55").Value = 4573
$ws.Range("M55").Value = -4258

$ws.Range("H58").Value = 4875.4116
$ws.Range("I58").Value = 3981.7856
$ws.Range("K58").Value = 3981.7856
$ws.Range("M58").Value = -3778.7856

$ws.Range("H64").Value = 75000
$ws.Range("J64").Value = 75000
$ws.Range("L64").Value = 75000
$ws.Range("N64").Value = -75496

$ws.Range("H67").Value = 75000
$ws.Range("J67").Value = 75000
$ws.Range("L67").Value = 75000
$ws.Range("N67").Value = -76716

$ws.Range("H68").Value = 41632.668
$ws.Range("J68").Value = 41632.668
$ws.Range("L68").Value = 41632.668
$ws.Range("N68").Value = -43130.668

$ws.Range("H71").Value = 41632.668
$ws.Range("J71").Value = 41632.668
$ws.Range("L71").Value = 124898.004
$ws.Range("N71").Value = -132386.004

$ws.Range("H99").Value = 5159.222
$ws.Range("I99").Value = 4429.7334
$ws.Range("J99").Value = 8806.666999999999
$ws.Range("K99").Value = 4429.7334
$ws.Range("L99").Value = 8806.666999999999
$ws.Range("M99").Value = -2931.7334
$ws.Range("N99").Value = -11802.667

$ws.Range("H126").Value = 5159.222
$ws.Range("I126").Value = 4429.7334
$ws.Range("J126").Value = 8806.666999999999
$ws.Range("K126").Value = 13289.2002
$ws.Range("L126").Value = 26420.001
$ws.Range("M126").Value = -10819.2002
$ws.Range("N126").Value = -31360.001

$ws.Range("H132").Value = 7505.3213
$ws.Range("I132").Value = 5461.364
$ws.Range("K132").Value = 16384.092
$ws.Range("M132").Value = -13854.092

$ws.Range("H134").Value = 2027.7222
$ws.Range("I134").Value = 1966.6
$ws.Range("J134").Value = 2333.3333
$ws.Range("K134").Value = 5899.799999999999
$ws.Range("L134").Value = 6999.999899999999
$ws.Range("M134").Value = -3364.799999999999
$ws.Range("N134").Value = -12069.9999

$ws.Range("H136").Value = 4875.4116
$ws.Range("I136").Value = 3981.7856
$ws.Range("K136").Value = 11945.3568
$ws.Range("M136").Value = -9395.356800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 285714460
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H121").Value = 1549.25
$ws.Range("I121").Value = 99.5
$ws.Range("J121").Value = 2999
$ws.Range("K121").Value = 298.5
$ws.Range("L121").Value = 8997
$ws.Range("M121").Value = 1011.5
$ws.Range("N121").Value = -11617

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492

$ws.Range("H97").Value = 1071.9
$ws.Range("I97").Value = 1117.8572
$ws.Range("K97").Value = 1117.8572
$ws.Range("M97").Value = -621.8571999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2283.2778
$ws.Range("I46").Value = 1512.375
$ws.Range("J46").Value = 2900
$ws.Range("K46").Value = 1512.375
$ws.Range("L46").Value = 2900
$ws.Range("M46").Value = -1324.375
$ws.Range("N46").Value = -3276

$ws.Range("H61").Value = 1919.6316
$ws.Range("I61").Value = 2030.3077
$ws.Range("J61").Value = 1679.8334
$ws.Range("K61").Value = 2030.3077
$ws.Range("L61").Value = 1679.8334
$ws.Range("M61").Value = -1828.3077
$ws.Range("N61").Value = -2083.8334

$ws.Range("H93").Value = 1287.9412
$ws.Range("I93").Value = 1181.5454
$ws.Range("J93").Value = 1483
$ws.Range("K93").Value = 1181.5454
$ws.Range("L93").Value = 1483
$ws.Range("M93").Value = 66.45460000000003
$ws.Range("N93").Value = -3979

$ws.Range("H113").Value = 1919.6316
$ws.Range("I113").Value = 2030.3077
$ws.Range("J113").Value = 1679.8334
$ws.Range("K113").Value = 2030.3077
$ws.Range("L113").Value = 1679.8334
$ws.Range("M113").Value = 139.6922999999999
$ws.Range("N113").Value = -6019.8334

$ws.Range("H116").Value = 252250
$ws.Range("J116").Value = 252250
$ws.Range("L116").Value = 252250
$ws.Range("N116").Value = -261428

$ws.Range("H132").Value = 2930.1667
$ws.Range("I132").Value = 3118.2
$ws.Range("J132").Value = 1990
$ws.Range("K132").Value = 9354.599999999999
$ws.Range("L132").Value = 5970
$ws.Range("M132").Value = -6824.599999999999
$ws.Range("N132").Value = -11030

$ws.Range("H136").Value = 2826.4546
$ws.Range("J136").Value = 3631.6667
$ws.Range("L136").Value = 10895.0001
$ws.Range("N136").Value = -15995.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 370.25
$ws.Range("J23").Value = 440.5
$ws.Range("L23").Value = 440.5
$ws.Range("N23").Value = -898.5

$ws.Range("H46").Value = 45000
$ws.Range("J46").Value = 45000
$ws.Range("L46").Value = 45000
$ws.Range("N46").Value = -45462

$ws.Range("H100").Value = 3670564.5
$ws.Range("I100").Value = 11617074
$ws.Range("J100").Value = 2944.923
$ws.Range("K100").Value = 23234148
$ws.Range("L100").Value = 5889.846
$ws.Range("M100").Value = -23233607
$ws.Range("N100").Value = -6971.846

$ws.Range("H107").Value = 485.8889
$ws.Range("J107").Value = 725
$ws.Range("L107").Value = 2175
$ws.Range("N107").Value = -6015

$ws.Range("H113").Value = 410.5
$ws.Range("I113").Value = 439.14285
$ws.Range("K113").Value = 1317.42855
$ws.Range("M113").Value = 852.5714499999999

$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 135000
$ws.Range("N134").Value = -140070
